$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-12-13 Saturday" "2025-12-14 Sunday"

Replace-Text "482÷9=" "244÷6="
Replace-Text "322÷4=" "392÷6="
Replace-Text "449÷3=" "595÷9="
Replace-Text "218÷4=" "280÷7="
Replace-Text "873÷5=" "891÷8="

Replace-Text "884÷4=" "391÷8="
Replace-Text "732÷6=" "349÷5="
Replace-Text "827÷7=" "301÷6="
Replace-Text "242÷7=" "459÷4="
Replace-Text "476÷7=" "369÷7="

Replace-Text "666÷7=" "326÷4="
Replace-Text "208÷3=" "113÷6="
Replace-Text "587÷2=" "503÷4="
Replace-Text "241÷8=" "543÷7="
Replace-Text "739÷7=" "514÷5="

Replace-Text "215÷3=" "779÷2="
Replace-Text "226÷5=" "542÷3="
Replace-Text "633÷3=" "905÷4="
Replace-Text "403÷3=" "729÷5="
Replace-Text "388÷4=" "413÷8="

Replace-Text "641÷4=" "398÷2="
Replace-Text "621÷5=" "594÷2="
Replace-Text "930÷5=" "352÷4="
Replace-Text "590÷9=" "833÷8="
Replace-Text "615÷3=" "808÷4="
